$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wsTsv = $wb.Worksheets.Item(1)
$wsTsv.Name = "Export as TSV"

$wsX = $wb.Worksheets.Item(5)
$wsX.Name = "ablation_dista...s_x_units list"

$wsY = $wb.Worksheets.Item(6)
$wsY.Name = "ablation_dista...s_y_units list"

# --- Freeze header row on the main sheet ---
$wsTsv.Activate()
$wsTsv.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Update data validations: error messages/titles, and list formulas that reference the renamed sheets ---

$wsTsv.Range("I2:I1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("I2:I1048576").Validation.ErrorMessage = "Value must be one of: mass_spectrometry_imaging."

$wsTsv.Range("J2:J1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("J2:J1048576").Validation.ErrorMessage = "Value must be one of: Imaging Mass Cytometry."

$wsTsv.Range("K2:K1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("K2:K1048576").Validation.ErrorMessage = "Value must be one of: protein."

$wsTsv.Range("L2:L1048576").Validation.ErrorTitle = "Not a boolean"
$wsTsv.Range("L2:L1048576").Validation.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$wsTsv.Range("S2:S1048576").Validation.ErrorTitle = "Not an integer"
$wsTsv.Range("S2:S1048576").Validation.ErrorMessage = "The values in this column must be integers."

$wsTsv.Range("T2:T1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("T2:T1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$wsTsv.Range("U2:U1048576").Validation.Formula1 = "='ablation_dista...s_x_units list'!`$A`$1:`$A`$2"
$wsTsv.Range("U2:U1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("U2:U1048576").Validation.ErrorMessage = "Value must be one of: um / nm."

$wsTsv.Range("V2:V1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("V2:V1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$wsTsv.Range("W2:W1048576").Validation.Formula1 = "='ablation_dista...s_y_units list'!`$A`$1:`$A`$2"
$wsTsv.Range("W2:W1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("W2:W1048576").Validation.ErrorMessage = "Value must be one of: um / nm."

$wsTsv.Range("X2:X1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("X2:X1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$wsTsv.Range("Y2:Y1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("Y2:Y1048576").Validation.ErrorMessage = "Value must be one of: Hz."

$wsTsv.Range("AA2:AA1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("AA2:AA1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$wsTsv.Range("AC2:AC1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("AC2:AC1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$wsTsv.Range("AE2:AE1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("AE2:AE1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$wsTsv.Range("AF2:AF1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("AF2:AF1048576").Validation.ErrorMessage = "Value must be one of: um."

$wsTsv.Range("AG2:AG1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("AG2:AG1048576").Validation.ErrorMessage = "The values in this column must be numbers."

$wsTsv.Range("AH2:AH1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("AH2:AH1048576").Validation.ErrorMessage = "Value must be one of: um."

$wsTsv.Range("AI2:AI1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("AI2:AI1048576").Validation.ErrorMessage = "Value must be one of: float / integer / string."

$wsTsv.Range("AJ2:AJ1048576").Validation.ErrorTitle = "Value must come from list"
$wsTsv.Range("AJ2:AJ1048576").Validation.ErrorMessage = "Value must be one of: dual count / pulse count / intensity value."

$wsTsv.Range("AL2:AL1048576").Validation.ErrorTitle = "Not a number"
$wsTsv.Range("AL2:AL1048576").Validation.ErrorMessage = "The values in this column must be numbers."
